$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)

# Resize / reposition the title placeholder (EMU 736979,628721,10686197,971479
# expressed in points; literals nudged slightly so the float32 COM marshalling
# truncates to the exact target EMU values).
$sh.Left   = 58.02984431968504
$sh.Top    = 49.5055905511811
$sh.Width  = 841.4328346456693
$sh.Height = 76.4944115488189

# Update the title text.
$sh.TextFrame.TextRange.Text = "mkarta.uz a medical record web site"

# Turn on "Shrink text on overflow" (normAutofit) for the text frame.
$sh.TextFrame.AutoSize = 2
